# Qualifier 1 completed with predictions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the actual scores for row 66 (Qualifier 1, MI vs DC) now that the match is complete ---
$ws.Range("E66").Value = 60
$ws.Range("H66").Value = 80
$ws.Range("K66").Value = 20
$ws.Range("N66").Value = 40
$ws.Range("Q66").Value = 100
$ws.Range("T66").Value = 0

# --- Coins earned from correct predictions (Qualifier 1 column) for the 5 existing players ---
$ws.Range("C81").Value = 10
$ws.Range("C82").Value = 7
$ws.Range("C83").Value = 0
$ws.Range("C84").Value = 10
$ws.Range("C85").Value = 5

# G81 becomes a SUM formula (previously a hard-coded count)
$ws.Range("G81").Formula = "=SUM(C81:F81)"

# G82:G85 get the same SUM formula, entered as one shared fill
$ws.Range("G82:G85").Formula = "=SUM(C82:F82)"

# --- Insert a new row for Balaji (6th player) below Sundar, before the grand-total row ---
$ws.Rows("86").Insert()

# Pull formatting for the new row from the row above (Sundar), then overwrite the content
$ws.Range("B85:G85").Copy()
$ws.Range("B86").PasteSpecial(-4122)
$ws.Range("L85:R85").Copy()
$ws.Range("L86").PasteSpecial(-4122)
$ws.Rows("86").RowHeight = 21

$ws.Range("B86").Value = "Balaji"
$ws.Range("C86").Value = 0
$ws.Range("G86").Formula = "=SUM(C86:F86)"

$ws.Range("L86").Value = "Balaji"
$ws.Range("M86").Formula = "=T73"
$ws.Range("N86").Value = 0
$ws.Range("O86").Formula = "=G86"
$ws.Range("P86").Formula = '=(-SUM($N$81:$N$85)/SUM($O$81:$O$85))*O86'
$ws.Range("Q86").Formula = "=SUM(M86,N86,P86)"
$ws.Range("R86").Value = "Balaji"

# Conditional formatting for the new row, matching the existing M81:M85 / Q81:Q85 rules
# (red for < 0, green for = 0, green for > 0)
$mLess = $ws.Range("M86").FormatConditions.Add(1, 6, "0")
$mLess.Interior.Color = 13551615
$mLess.Font.Color = 393372
$mEqual = $ws.Range("M86").FormatConditions.Add(1, 3, "0")
$mEqual.Interior.Color = 13561798
$mEqual.Font.Color = 24832
$mGreater = $ws.Range("M86").FormatConditions.Add(1, 5, "0")
$mGreater.Interior.Color = 13561798
$mGreater.Font.Color = 24832

$qLess = $ws.Range("Q86").FormatConditions.Add(1, 6, "0")
$qLess.Interior.Color = 13551615
$qLess.Font.Color = 393372
$qEqual = $ws.Range("Q86").FormatConditions.Add(1, 3, "0")
$qEqual.Interior.Color = 13561798
$qEqual.Font.Color = 24832
$qGreater = $ws.Range("Q86").FormatConditions.Add(1, 5, "0")
$qGreater.Interior.Color = 13561798
$qGreater.Font.Color = 24832

# --- Fix up the final totals formula (old row 86, now shifted to row 87) ---
$ws.Range("Q87").Formula = "=SUM(Q81:Q85,T73)"

# --- Keep selection in sync with the shifted totals cell ---
[void]$ws.Range("Q87").Select()

$wb.Save()
